$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped by GitHub Actions
# (also includes the Aave / NEARProtocol row swap at rows 47-48).
# Force the Price/Coin/Link columns to stay text so values such as "0.9980"
# or "34.03" are not reinterpreted as numbers and lose trailing zeros.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "B47", "C47", "D47", "E47", "B48", "C48", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellAddr in $cells) { $ws.Range($cellAddr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.441.53"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.896.67"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "237.55"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "0.4851"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "0.2906"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").Value = "0.06621"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "1.910.17"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "16.99"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "0.07329"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "5.181"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "88.01"
$ws.Range("D15").Value = "0.6633"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "30.424.74"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "13.47"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "0.000007791"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "0.9968"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "5.441"
$ws.Range("E20").Value = "  +4.38%  "
$ws.Range("D21").Value = "2.125.73"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "0.9980"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "195.26"
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("D24").Value = "6.201"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "9.354"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").Value = "165.22"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "18.22"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("D28").Value = "1.947"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "1.450"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "4.312"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "0.09174"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "4.056"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "0.05091"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "1.160"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").Value = "0.7283"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "0.01791"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "2.652"
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("D39").Value = "0.9216"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "2.091"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "106.29"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "0.4327"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").Value = "5.869"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").Value = "1.004"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "7.551"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "0.1324"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "1.579"
$ws.Range("E47").Value = "  +9.43%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "65.19"
$ws.Range("E48").Value = "  -10.42%  "
$ws.Range("D49").Value = "8.990"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "34.03"
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").Value = "0.05762"
$ws.Range("E51").Value = "  -3.15%  "
